# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# Offense sheet ("OFF") - Road row (row 3) totals updated to include Week 17 data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 247
$wsOff.Range("C3").Value = 174
$wsOff.Range("D3").Value = 68
$wsOff.Range("E3").Value = 32
$wsOff.Range("F3").Value = 6
$wsOff.Range("G3").Value = 7

# Defense sheet ("DEF") - Road row (row 3) totals updated to include Week 17 data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 283
$wsDef.Range("C3").Value = 214
$wsDef.Range("D3").Value = 66
$wsDef.Range("E3").Value = 33
$wsDef.Range("F3").Value = 5
$wsDef.Range("G3").Value = 2
